$wb = $excel.ActiveWorkbook

# 1. Remove the GRADE and STU_GRADE design sheets - no longer part of the schema doc.
$wb.Worksheets.Item("年级表(GRADE)").Delete()
$wb.Worksheets.Item("学生_年级表（STU_GRADE)").Delete()

# 2. USER sheet: just a cursor/selection move, no data change.
$wsUser = $wb.Worksheets.Item("用户表设计(USER)")
$wsUser.Range("A7").Select()

# 3. EXPRESSION sheet: add the new e_tip (做题提示) field.
$wsExpr = $wb.Worksheets.Item("表达式(EXPRESSION)")
$wsExpr.Range("A6").Value = "e_tip"
$wsExpr.Range("B6").Value = "做题提示"
$wsExpr.Range("C6").Value = "VARCHAR"
$wsExpr.Range("D6").Value = "长度<=100,可以为空"

# 4. record sheet: r_time's type is now VARCHAR(25) instead of TIME, plus two new
#    fields (r_num, r_rank) describing attempt count / difficulty rank.
$wsRecord = $wb.Worksheets.Item("做题记录表(record)")
$wsRecord.Range("C6").Value = "VARCHAR(25)"

$wsRecord.Range("A7").Value = "r_num"
$wsRecord.Range("B7").Value = "用户是第几次做题"
$wsRecord.Range("C7").Value = "INT "
$wsRecord.Range("D7").Value = "不为空"

$wsRecord.Range("A8").Value = "r_rank"
$wsRecord.Range("B8").Value = "题目难度等级"
$wsRecord.Range("C8").Value = "INT"
$wsRecord.Range("D8").Value = "不为空 "
$wsRecord.Range("E8").Value = "分为五个等级 ， 分别是 1，2，3，4，5"

# Column widths for the newly meaningful B/C columns (closest reachable values).
$wsRecord.Columns.Item(2).ColumnWidth = 14.4
$wsRecord.Columns.Item(3).ColumnWidth = 12.0

# record sheet becomes the active tab again, parked at E9 this time.
$wsRecord.Range("E9").Select()
